# Insert a new data row at row 168 (pushes existing rows 168-254 down to 169-255)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(168).EntireRow.Insert()

$ws.Range("A168").Value = 5
$ws.Range("B168").Value = "Macroferia Regional de Talca"
$ws.Range("C168").Value = "Maule"
$ws.Range("D168").Value = 44719
$ws.Range("E168").Value = 7
$ws.Range("F168").Value = 100112045
$ws.Range("G168").Value = "Zapallo"
$ws.Range("H168").Value = "Camote"
$ws.Range("I168").Value = "1a (guarda)"
$ws.Range("J168").Value = 900
$ws.Range("K168").Value = 350
$ws.Range("L168").Value = 350
$ws.Range("M168").Value = 350
$ws.Range("N168").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O168").Value = "Región del Maule"
$ws.Range("P168").Value = 350
$ws.Range("Q168").Value = 1
$ws.Range("R168").Value = "Hortaliza"

# Preserve the same number style/format as the other date cells in column D
$ws.Range("D168").NumberFormat = $ws.Range("D169").NumberFormat
